$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose values look numeric must be forced to Text format so that
# Excel does not silently convert the stored string into a floating point
# number (the source data keeps its original textual formatting, e.g.
# trailing zeros, fixed decimal places, thousands separators, etc).

$ws.Range("D2").Value = '19.986.97'
$ws.Range("E2").Value = '  -7.05%  '

$ws.Range("D3").Value = '1.413.58'
$ws.Range("E3").Value = '  -7.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9987'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '275.26'
$ws.Range("E6").Value = '  -4.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3670'
$ws.Range("E7").Value = '  -5.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3117'
$ws.Range("E8").Value = '  -1.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.87'
$ws.Range("E9").Value = '  -6.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.040'
$ws.Range("E10").Value = '  -2.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06506'
$ws.Range("E11").Value = '  -8.99%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.502'
$ws.Range("E13").Value = '  -3.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.76'
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.207'
$ws.Range("E15").Value = '  -4.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001023'

$ws.Range("D17").Value = '1.410.13'
$ws.Range("E17").Value = '  -8.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05697'
$ws.Range("E18").Value = '  -13.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.0000'
$ws.Range("E19").Value = '  -0.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.06'
$ws.Range("E20").Value = '  -14.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.647'
$ws.Range("E21").Value = '  -7.36%  '

$ws.Range("D22").Value = '14.74'
$ws.Range("E22").Value = '  -4.15%  '

$ws.Range("D23").Value = '11.09'
$ws.Range("E23").Value = '  +2.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.273'
$ws.Range("E24").Value = '  -4.09%  '

$ws.Range("D25").Value = '20.016.87'
$ws.Range("E25").Value = '  -6.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.277'
$ws.Range("E26").Value = '  -3.98%  '

$ws.Range("D27").Value = '133.59'
$ws.Range("E27").Value = '  -10.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.13'
$ws.Range("E28").Value = '  -6.54%  '

$ws.Range("D29").Value = '1.570.86'
$ws.Range("E29").Value = '  -8.06%  '

$ws.Range("D30").Value = '109.86'
$ws.Range("E30").Value = '  -5.68%  '

$ws.Range("D31").Value = '3.972'
$ws.Range("E31").Value = '  -17.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.320'
$ws.Range("E32").Value = '  -12.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8284'
$ws.Range("E33").Value = '  -13.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07694'
$ws.Range("E34").Value = '  -3.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.405'
$ws.Range("E35").Value = '  -0.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.483'
$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05919'
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.919'
$ws.Range("E38").Value = '  -4.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.000'
$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02075'
$ws.Range("E40").Value = '  -5.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.54'
$ws.Range("E41").Value = '  -6.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1907'
$ws.Range("E42").Value = '  -5.53%  '

$ws.Range("D43").Value = '1.099'
$ws.Range("E43").Value = '  -6.85%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.40'
$ws.Range("E44").Value = '  -5.18%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5327'
$ws.Range("E45").Value = '  -7.25%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.531'
$ws.Range("E46").Value = '  -5.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5221'
$ws.Range("E47").Value = '  -5.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '116.57'
$ws.Range("E48").Value = '  +1.01%  '

$ws.Range("D49").Value = '1.774'
$ws.Range("E49").Value = '  -6.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.041'
$ws.Range("E50").Value = '  -10.03%  '

$ws.Range("E51").Value = '  -0.13%  '
